# Update generated output values (gh-pages regeneration) for
# sheets "展览" (sheet1) and "全部类型" (sheet4).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 22
$ws1.Range("F6").Value = 16
$ws1.Range("F7").Value = 560
$ws1.Range("F8").Value = 7876
$ws1.Range("F9").Value = 750
$ws1.Range("F10").Value = 221
$ws1.Range("F11").Value = 1094
$ws1.Range("F12").Value = 740
$ws1.Range("F15").Value = 195
$ws1.Range("F16").Value = 39
$ws1.Range("I17").Value = "//i0.hdslb.com/bfs/openplatform/202405/MiqOsFGU1715224161257.jpeg"
$ws1.Range("F18").Value = 811

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 22
$ws4.Range("F6").Value = 16
$ws4.Range("F8").Value = 560
$ws4.Range("F9").Value = 7876
$ws4.Range("F10").Value = 750
$ws4.Range("F11").Value = 221
$ws4.Range("F12").Value = 1094
$ws4.Range("F13").Value = 740
$ws4.Range("F16").Value = 195
$ws4.Range("F17").Value = 39
$ws4.Range("I18").Value = "//i0.hdslb.com/bfs/openplatform/202405/MiqOsFGU1715224161257.jpeg"
$ws4.Range("F19").Value = 811
